# Re-order the storage options listed in the "Azure PaaS noSQL Database"
# rectangle on slide 1:
#   before: (Table storage, DocumentDB, HDInsight Hbase, ...)
#   after : (DocumentDB, HDInsight Hbase, Table storage, ...)
#
# The run-splitting in the target OOXML shows that PowerPoint reused the
# five pre-existing runs (keeping each run's own rPr, i.e. its dirty/err
# spell-check flags) and only split two of them in two to make room for
# the extra word-boundaries introduced by the reorder. We reproduce that
# by editing each run's text *in place* (never touching text that spans
# more than one original run at a time), relying on the engine to split
# a run cleanly whenever we set the .Text of a narrower sub-range.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$ellipsis = [char]0x2026

# --- Run "(Table " -> "(DocumentDB, HDInsight " -------------------------
$full = $tr.Text
$idx = $full.IndexOf("(Table ")
$start = $idx + 1
$run = $tr.Characters($start, 7)
$run.Text = "(DocumentDB, HDInsight "

# --- Run "storage" -> "Hbase" ---------------------------------------------
$full = $tr.Text
$idx = $full.IndexOf("storage")
$start = $idx + 1
$run = $tr.Characters($start, 7)
$run.Text = "Hbase"

# --- Run ", DocumentDB, HDInsight " -> split into ", " + "Table " --------
$full = $tr.Text
$idx = $full.IndexOf(", DocumentDB, HDInsight ")
$start = $idx + 1
$run = $tr.Characters($start, 2)
$run.Text = ", "
$start2 = $start + 2
$run2 = $tr.Characters($start2, 22)
$run2.Text = "Table "

# --- Run "Hbase" (now directly after "Table ") -> "storage" --------------
$full = $tr.Text
$idx = $full.IndexOf("Table Hbase")
$start = $idx + 1 + ("Table ").Length
$run = $tr.Characters($start, 5)
$run.Text = "storage"

# --- Run ", ...)" (last 4 chars) -> split into ", " + "...)" -------------
$len = $tr.Length
$start = $len - 3
$run = $tr.Characters($start, 2)
$run.Text = ", "
$start2 = $start + 2
$run2 = $tr.Characters($start2, 2)
$run2.Text = $ellipsis + ")"

Write-Host "Final text: $($tr.Text)"
